# Auto-generated: update FFXIV Leve profit-tracking sheets with refreshed
# market-board figures (columns H-N) as produced by the scheduled data-pull runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8308.733
$ws.Range("I32").Value = 7944
$ws.Range("J32").Value = 8441.362999999999
$ws.Range("K32").Value = 7944
$ws.Range("L32").Value = 8441.362999999999
$ws.Range("M32").Value = -7618
$ws.Range("N32").Value = -9093.362999999999
$ws.Range("H92").Value = 44721.875
$ws.Range("I92").Value = 17624.69
$ws.Range("J92").Value = 306661.34
$ws.Range("K92").Value = 17624.69
$ws.Range("L92").Value = 306661.34
$ws.Range("M92").Value = -16376.69
$ws.Range("N92").Value = -309157.34
$ws.Range("H107").Value = 1326.8667
$ws.Range("J107").Value = 1377
$ws.Range("L107").Value = 1377
$ws.Range("N107").Value = -5217
$ws.Range("H115").Value = 1128.7273
$ws.Range("I115").Value = 951.6
$ws.Range("J115").Value = 2900
$ws.Range("K115").Value = 2854.8
$ws.Range("L115").Value = 8700
$ws.Range("M115").Value = -1287.8
$ws.Range("N115").Value = -11834
$ws.Range("H132").Value = 1829.4906
$ws.Range("I132").Value = 1720.4584
$ws.Range("K132").Value = 5161.3752
$ws.Range("M132").Value = -2631.3752
$ws.Range("H139").Value = 64138.332
$ws.Range("J139").Value = 64138.332
$ws.Range("L139").Value = 64138.332
$ws.Range("N139").Value = -74418.33199999999
$ws.Range("H140").Value = 99342.28999999999
$ws.Range("J140").Value = 99342.28999999999
$ws.Range("L140").Value = 99342.28999999999
$ws.Range("N140").Value = -109702.29

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8050.86
$ws.Range("I32").Value = 5845.225
$ws.Range("K32").Value = 5845.225
$ws.Range("M32").Value = -5558.225
$ws.Range("H33").Value = 19999
$ws.Range("I33").Value = 19999
$ws.Range("K33").Value = 19999
$ws.Range("M33").Value = -19670
$ws.Range("H37").Value = 21196.4
$ws.Range("J37").Value = 49999
$ws.Range("L37").Value = 49999
$ws.Range("N37").Value = -50545
$ws.Range("H45").Value = 659.6326
$ws.Range("I45").Value = 618.9318
$ws.Range("J45").Value = 1017.8
$ws.Range("K45").Value = 618.9318
$ws.Range("L45").Value = 1017.8
$ws.Range("M45").Value = -241.9318
$ws.Range("N45").Value = -1771.8
$ws.Range("H74").Value = 2796.1667
$ws.Range("I74").Value = 2508.3333
$ws.Range("J74").Value = 3659.6667
$ws.Range("K74").Value = 2508.3333
$ws.Range("L74").Value = 3659.6667
$ws.Range("M74").Value = -1634.3333
$ws.Range("N74").Value = -5407.6667
$ws.Range("H77").Value = 2796.1667
$ws.Range("I77").Value = 2508.3333
$ws.Range("J77").Value = 3659.6667
$ws.Range("K77").Value = 12541.6665
$ws.Range("L77").Value = 18298.3335
$ws.Range("M77").Value = -8173.666499999999
$ws.Range("N77").Value = -27034.3335
$ws.Range("H132").Value = 11300.156
$ws.Range("I132").Value = 11857.714
$ws.Range("J132").Value = 7397.25
$ws.Range("K132").Value = 35573.142
$ws.Range("L132").Value = 22191.75
$ws.Range("M132").Value = -33043.142
$ws.Range("N132").Value = -27251.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2252.516
$ws.Range("J134").Value = 2500
$ws.Range("L134").Value = 7500
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 172
$ws.Range("I7").Value = 115.71429
$ws.Range("K7").Value = 115.71429
$ws.Range("M7").Value = -2.714290000000005
$ws.Range("H31").Value = 2989.8667
$ws.Range("I31").Value = 2769.4707
$ws.Range("K31").Value = 2769.4707
$ws.Range("M31").Value = -2474.4707
$ws.Range("H34").Value = 2989.8667
$ws.Range("I34").Value = 2769.4707
$ws.Range("K34").Value = 2769.4707
$ws.Range("M34").Value = -2567.4707
$ws.Range("H52").Value = 94700
$ws.Range("J52").Value = 99400
$ws.Range("L52").Value = 99400
$ws.Range("N52").Value = -99988
$ws.Range("H105").Value = 633.3182
$ws.Range("I105").Value = 636.75
$ws.Range("K105").Value = 636.75
$ws.Range("M105").Value = 1110.25
$ws.Range("H129").Value = 45280
$ws.Range("J129").Value = 45280
$ws.Range("L129").Value = 45280
$ws.Range("N129").Value = -55280
$ws.Range("I132").Value = 3658.2144
$ws.Range("J132").Value = 27496
$ws.Range("K132").Value = 10974.6432
$ws.Range("L132").Value = 82488
$ws.Range("M132").Value = -8444.643199999999
$ws.Range("N132").Value = -87548
$ws.Range("H135").Value = 67936.25
$ws.Range("J135").Value = 67936.25
$ws.Range("L135").Value = 67936.25
$ws.Range("N135").Value = -78076.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H124").Value = 4239.6665
$ws.Range("I124").Value = 2924.3333
$ws.Range("J124").Value = 5555
$ws.Range("K124").Value = 8772.999899999999
$ws.Range("L124").Value = 16665
$ws.Range("M124").Value = -3862.999899999999
$ws.Range("N124").Value = -26485
$ws.Range("H129").Value = 1274
$ws.Range("I129").Value = 821.5
$ws.Range("J129").Value = 1952.75
$ws.Range("K129").Value = 2464.5
$ws.Range("L129").Value = 5858.25
$ws.Range("M129").Value = 2535.5
$ws.Range("N129").Value = -15858.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 30500
$ws.Range("J38").Value = 30500
$ws.Range("L38").Value = 30500
$ws.Range("N38").Value = -31426
$ws.Range("H80").Value = 9651.6
$ws.Range("I80").Value = 4174.25
$ws.Range("J80").Value = 14707.615
$ws.Range("K80").Value = 4174.25
$ws.Range("L80").Value = 14707.615
$ws.Range("M80").Value = -3176.25
$ws.Range("N80").Value = -16703.615
$ws.Range("H83").Value = 9651.6
$ws.Range("I83").Value = 4174.25
$ws.Range("J83").Value = 14707.615
$ws.Range("K83").Value = 20871.25
$ws.Range("L83").Value = 73538.075
$ws.Range("M83").Value = -15879.25
$ws.Range("N83").Value = -83522.075
$ws.Range("H126").Value = 85083.16
$ws.Range("J126").Value = 4919.4
$ws.Range("L126").Value = 14758.2
$ws.Range("N126").Value = -19698.2
$ws.Range("H132").Value = 4191.294
$ws.Range("I132").Value = 4403.8184
$ws.Range("J132").Value = 3801.6667
$ws.Range("K132").Value = 13211.4552
$ws.Range("L132").Value = 11405.0001
$ws.Range("M132").Value = -10681.4552
$ws.Range("N132").Value = -16465.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 45012.105
$ws.Range("J40").Value = 14747.9
$ws.Range("L40").Value = 14747.9
$ws.Range("N40").Value = -15019.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2873.4849
$ws.Range("I132").Value = 1944.7368
$ws.Range("K132").Value = 5834.2104
$ws.Range("M132").Value = -3304.2104

